$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap Candidate (column D) and Action Date (column F) values between row 4 and row 5
$d4 = $ws.Range("D4").Value2
$f4 = $ws.Range("F4").Value2
$d5 = $ws.Range("D5").Value2
$f5 = $ws.Range("F5").Value2

$ws.Range("D4").Value2 = $d5
$ws.Range("F4").Value2 = $f5
$ws.Range("D5").Value2 = $d4
$ws.Range("F5").Value2 = $f4
